$d = $word.ActiveDocument

# 1) Programa - Portuguese paragraph: turn the run-on "1 - ... 9 - ..." text
#    into one <w:t> per numbered item separated by manual line breaks (^l -> <w:br/>).
$found1 = $d.Content.Find.Execute('1 - Introdução ao Laboratório: Noções Elementares de Segurança; Equipamentos Básicos de Laboratório; Equipamentos de Proteção Individual.2 - Pesos e medidas (Tratamento de dados experimentais): Cuidados Gerais com Balanças; Técnicas de Determinação de massa; Exatidão e precisão; Unidades; Algarismos Significativos; Propagação de Erros.3 - Técnicas de Separação de Misturas: Filtração simples; Filtração a vácuo e Decantação.4 - Fenômenos físicos: Construção do Diagrama da mudança do estado físico da água.5 - Miscibilidade e solubilidade: Influência das forças intermoleculares na miscibilidade de líquidos.6 - Reações químicas: Aspectos qualitativos.7 - Soluções: Preparo e padronização de soluções.8 - Titrimetria: Realização de Titulações Ácido-Base; Retrotitulação.9 - Equilíbrio Químico - Preparo de Solução Tampão.', $false, $false, $false, $false, $false, $true, 1, $false, '1 - Introdução ao Laboratório: Noções Elementares de Segurança; Equipamentos Básicos de Laboratório; Equipamentos de Proteção Individual.^l2 - Pesos e medidas (Tratamento de dados experimentais): Cuidados Gerais com Balanças; Técnicas de Determinação de massa; Exatidão e precisão; Unidades; Algarismos Significativos; Propagação de Erros.^l3 - Técnicas de Separação de Misturas: Filtração simples; Filtração a vácuo e Decantação.^l4 - Fenômenos físicos: Construção do Diagrama da mudança do estado físico da água.^l5 - Miscibilidade e solubilidade: Influência das forças intermoleculares na miscibilidade de líquidos.^l6 - Reações químicas: Aspectos qualitativos.^l7 - Soluções: Preparo e padronização de soluções.^l8 - Titrimetria: Realização de Titulações Ácido-Base; Retrotitulação.^l9 - Equilíbrio Químico - Preparo de Solução Tampão.', 2)
Write-Output "Programa (PT) replaced: $found1"

# 2) Programa - English (italic) paragraph: same split, single line breaks.
$found2 = $d.Content.Find.Execute('1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment.2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation.3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation.4 - Physical phenomena: Water state changes.5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. 6 - Chemical reactions: Qualitative aspects.7 - Solutions: Preparation and standardization of solutions.8 - Titrimetry: Acid-Base Titrations and return-titration.9 - Chemical equilibrium: Buffer solution.', $false, $false, $false, $false, $false, $true, 1, $false, '1 - Introduction to the Chemistry Laboratory: Elementary notion of security, Laboratory basic equipment; Individual protection equipment.^l2 - Weights and measures (experimental data treatment): General care with scales, Determination of mass techniques. Accuracy and precision, units, significant digits and error propagation.^l3 - Methods for separating mixtures: Simple filtration; Vacuum filtration and Decantation.^l4 - Physical phenomena: Water state changes.^l5 - Miscibility and solubility: Intermolecular forces influence on the liquids miscibility. ^l6 - Chemical reactions: Qualitative aspects.^l7 - Solutions: Preparation and standardization of solutions.^l8 - Titrimetry: Acid-Base Titrations and return-titration.^l9 - Chemical equilibrium: Buffer solution.', 2)
Write-Output "Programa (EN) replaced: $found2"

# 3) Bibliografia paragraph: split each reference onto its own line, separated by a
#    blank line (two manual line breaks, ^l^l -> <w:br/><w:br/>).
$found3 = $d.Content.Find.Execute('ASSUMPÇÃO, R. M. V. ; MORITA, T. Manual de soluções reagentes e solventes: padronização, preparação, purificação. São Paulo: Editora Edgard Blucher, 1972.BACCAN, N.; ANDRADE, J. C. O. ; GODINHO, E. S.; BARONE, J. S. Química analítica quantitativa elementar. 2.ed. São Paulo: Edgard Blucher, 1995.BRADY, J; HUMISTON, G. E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1986.BROWN, T. E et al. Química a Ciência Central. 9 ed. São Paulo. Pearson Prentice Hall, 2005-2007.CONSTANTINO, M.G; SILVA, G. V. J. da; DONATE P. M. Fundamentos de química experimental, São Paulo : EDUSP, 2004.MAHAN, B. M.; MYERS, R. J. Química um curso universitário. São Paulo: Ed. Edgard Blucher Ltda, 1993.SILVA, R. R.; BOCCHI, N. ; ROCHA FILHO, R. C. Introdução a química experimental. São Paulo: McGraw-Hill, 1990.', $false, $false, $false, $false, $false, $true, 1, $false, 'ASSUMPÇÃO, R. M. V. ; MORITA, T. Manual de soluções reagentes e solventes: padronização, preparação, purificação. São Paulo: Editora Edgard Blucher, 1972.^l^lBACCAN, N.; ANDRADE, J. C. O. ; GODINHO, E. S.; BARONE, J. S. Química analítica quantitativa elementar. 2.ed. São Paulo: Edgard Blucher, 1995.^l^lBRADY, J; HUMISTON, G. E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1986.^l^lBROWN, T. E et al. Química a Ciência Central. 9 ed. São Paulo. Pearson Prentice Hall, 2005-2007.^l^lCONSTANTINO, M.G; SILVA, G. V. J. da; DONATE P. M. Fundamentos de química experimental, São Paulo : EDUSP, 2004.^l^lMAHAN, B. M.; MYERS, R. J. Química um curso universitário. São Paulo: Ed. Edgard Blucher Ltda, 1993.^l^lSILVA, R. R.; BOCCHI, N. ; ROCHA FILHO, R. C. Introdução a química experimental. São Paulo: McGraw-Hill, 1990.', 2)
Write-Output "Bibliografia replaced: $found3"
